$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All D/E (and the B/C touched below) columns hold explicit text values in the
# source workbook (t="inlineStr"), e.g. "1.001" or "22.021.41" -- force text
# number format before assignment so COM does not auto-coerce them to numbers.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "22.042.07"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "1.555.45"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "286.91"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3762"
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3244"
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "41.41"
$ws.Range("E9").Value = "  -12.87%  "
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "1.127"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07306"
$ws.Range("E11").Value = "  -2.84%  "
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "19.58"
$ws.Range("E13").Value = "  -5.81%  "
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "5.711"
$ws.Range("E14").Value = "  -3.76%  "
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "6.864"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "1.552.74"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001082"
$ws.Range("E17").Value = "  -3.27%  "
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06647"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "85.21"
$ws.Range("E19").Value = "  -3.56%  "
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "6.445"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "16.00"
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "22.062.67"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "2.243"
$ws.Range("E25").Value = "  -6.28%  "
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "2.534"
$ws.Range("E26").Value = "  -3.49%  "
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "149.89"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = "18.87"
$ws.Range("E28").Value = "  -3.92%  "
$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D29").Value = "4.841"
$ws.Range("E29").Value = "  -2.31%  "
$ws.Range("D30:E30").NumberFormat = "@"
$ws.Range("D30").Value = "1.728.31"
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31:E31").NumberFormat = "@"
$ws.Range("D31").Value = "120.14"
$ws.Range("E31").Value = "  -4.12%  "
$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("D32").Value = "1.120"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("D33").Value = "5.934"
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("B34:E34").NumberFormat = "@"
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "1.659"
$ws.Range("E34").Value = "  -16.59%  "
$ws.Range("B35:E35").NumberFormat = "@"
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "9.270"
$ws.Range("E35").Value = "  -5.82%  "
$ws.Range("B36:E36").NumberFormat = "@"
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "0.08114"
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("D37").Value = "5.233"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02288"
$ws.Range("E38").Value = "  -6.98%  "
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06134"
$ws.Range("E39").Value = "  -4.13%  "
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2116"
$ws.Range("E40").Value = "  -5.25%  "
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "1.219"
$ws.Range("E41").Value = "  -7.42%  "
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "10.92"
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5949"
$ws.Range("E44").Value = "  -5.15%  "
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "13.55"
$ws.Range("E45").Value = "  -3.73%  "
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "3.727"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5752"
$ws.Range("E47").Value = "  -5.61%  "
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "1.946"
$ws.Range("E48").Value = "  -5.00%  "
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "120.04"
$ws.Range("E49").Value = "  -4.00%  "
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "1.156"
$ws.Range("E50").Value = "  -4.43%  "
$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06933"
$ws.Range("E51").Value = "  -3.78%  "
